# Natmi following Dr Hou advice
# Expand the Psen1 -> Notch3 sending/target cluster matrix from a partial
# (ECs/FAPs senders only) set to the full 4x4 ECs/FAPs/M2/sCs cross product,
# and refresh the NATMI specificity metrics for every pair.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Psen1"
$ws.Range("C2").Value = "Notch3"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 27.58598166666667
$ws.Range("H2").Value = 82.75794500000001
$ws.Range("I2").Value = 0.2704460545904799
$ws.Range("J2").Value = 0.2704460545904799
$ws.Range("K2").Value = 2
$ws.Range("L2").Value = 0.6666666666666666
$ws.Range("M2").Value = 3.684362666666666
$ws.Range("N2").Value = 11.053088
$ws.Range("O2").Value = 0.0397011572965827
$ws.Range("P2").Value = 0.03970115729658269
$ws.Range("Q2").Value = 101.6367609760178
$ws.Range("R2").Value = 914.73084878416
$ws.Range("S2").Value = 0.01073702135353683
$ws.Range("T2").Value = 0.01073702135353683

# Row 3
$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Psen1"
$ws.Range("C3").Value = "Notch3"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 27.58598166666667
$ws.Range("H3").Value = 82.75794500000001
$ws.Range("I3").Value = 0.2704460545904799
$ws.Range("J3").Value = 0.2704460545904799
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 5.038243666666667
$ws.Range("N3").Value = 15.114731
$ws.Range("O3").Value = 0.05429001496473517
$ws.Range("P3").Value = 0.05429001496473517
$ws.Range("Q3").Value = 138.9848974208661
$ws.Range("R3").Value = 1250.864076787795
$ws.Range("S3").Value = 0.01468252035087074
$ws.Range("T3").Value = 0.01468252035087074

# Row 4
$ws.Range("A4").Value = "ECs"
$ws.Range("B4").Value = "Psen1"
$ws.Range("C4").Value = "Notch3"
$ws.Range("D4").Value = "M2"
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 27.58598166666667
$ws.Range("H4").Value = 82.75794500000001
$ws.Range("I4").Value = 0.2704460545904799
$ws.Range("J4").Value = 0.2704460545904799
$ws.Range("K4").Value = 2
$ws.Range("L4").Value = 0.6666666666666666
$ws.Range("M4").Value = 0.09716666666666667
$ws.Range("N4").Value = 0.2915
$ws.Range("O4").Value = 0.001047027523164011
$ws.Range("P4").Value = 0.001047027523164011
$ws.Range("Q4").Value = 2.680437885277778
$ws.Range("R4").Value = 24.1239409675
$ws.Range("S4").Value = 0.0002831644626873492
$ws.Range("T4").Value = 0.0002831644626873491

# Row 5
$ws.Range("A5").Value = "ECs"
$ws.Range("B5").Value = "Psen1"
$ws.Range("C5").Value = "Notch3"
$ws.Range("D5").Value = "sCs"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 27.58598166666667
$ws.Range("H5").Value = 82.75794500000001
$ws.Range("I5").Value = 0.2704460545904799
$ws.Range("J5").Value = 0.2704460545904799
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 83.98262666666666
$ws.Range("N5").Value = 251.94788
$ws.Range("O5").Value = 0.9049618002155182
$ws.Range("P5").Value = 0.9049618002155182
$ws.Range("Q5").Value = 2316.743199545178
$ws.Range("R5").Value = 20850.6887959066
$ws.Range("S5").Value = 0.244743348423385
$ws.Range("T5").Value = 0.244743348423385

# Row 6
$ws.Range("A6").Value = "FAPs"
$ws.Range("B6").Value = "Psen1"
$ws.Range("C6").Value = "Notch3"
$ws.Range("D6").Value = "ECs"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 26.23504533333333
$ws.Range("H6").Value = 78.705136
$ws.Range("I6").Value = 0.2572018131577233
$ws.Range("J6").Value = 0.2572018131577233
$ws.Range("K6").Value = 2
$ws.Range("L6").Value = 0.6666666666666666
$ws.Range("M6").Value = 3.684362666666666
$ws.Range("N6").Value = 11.053088
$ws.Range("O6").Value = 0.0397011572965827
$ws.Range("P6").Value = 0.03970115729658269
$ws.Range("Q6").Value = 96.65942158444088
$ws.Range("R6").Value = 869.9347942599678
$ws.Range("S6").Value = 0.01021120964114105
$ws.Range("T6").Value = 0.01021120964114104

# Row 7
$ws.Range("A7").Value = "FAPs"
$ws.Range("B7").Value = "Psen1"
$ws.Range("C7").Value = "Notch3"
$ws.Range("D7").Value = "FAPs"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 26.23504533333333
$ws.Range("H7").Value = 78.705136
$ws.Range("I7").Value = 0.2572018131577233
$ws.Range("J7").Value = 0.2572018131577233
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 5.038243666666667
$ws.Range("N7").Value = 15.114731
$ws.Range("O7").Value = 0.05429001496473517
$ws.Range("P7").Value = 0.05429001496473517
$ws.Range("Q7").Value = 132.1785509953796
$ws.Range("R7").Value = 1189.606958958416
$ws.Range("S7").Value = 0.01396349028528982
$ws.Range("T7").Value = 0.01396349028528982

# Row 8
$ws.Range("A8").Value = "FAPs"
$ws.Range("B8").Value = "Psen1"
$ws.Range("C8").Value = "Notch3"
$ws.Range("D8").Value = "M2"
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 26.23504533333333
$ws.Range("H8").Value = 78.705136
$ws.Range("I8").Value = 0.2572018131577233
$ws.Range("J8").Value = 0.2572018131577233
$ws.Range("K8").Value = 2
$ws.Range("L8").Value = 0.6666666666666666
$ws.Range("M8").Value = 0.09716666666666667
$ws.Range("N8").Value = 0.2915
$ws.Range("O8").Value = 0.001047027523164011
$ws.Range("P8").Value = 0.001047027523164011
$ws.Range("Q8").Value = 2.549171904888889
$ws.Range("R8").Value = 22.942547144
$ws.Range("S8").Value = 0.0002692973773838239
$ws.Range("T8").Value = 0.0002692973773838238

# Row 9
$ws.Range("A9").Value = "FAPs"
$ws.Range("B9").Value = "Psen1"
$ws.Range("C9").Value = "Notch3"
$ws.Range("D9").Value = "sCs"
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 26.23504533333333
$ws.Range("H9").Value = 78.705136
$ws.Range("I9").Value = 0.2572018131577233
$ws.Range("J9").Value = 0.2572018131577233
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 83.98262666666666
$ws.Range("N9").Value = 251.94788
$ws.Range("O9").Value = 0.9049618002155182
$ws.Range("P9").Value = 0.9049618002155182
$ws.Range("Q9").Value = 2203.288017812409
$ws.Range("R9").Value = 19829.59216031168
$ws.Range("S9").Value = 0.2327578158539087
$ws.Range("T9").Value = 0.2327578158539086

# Row 10
$ws.Range("A10").Value = "M2"
$ws.Range("B10").Value = "Psen1"
$ws.Range("C10").Value = "Notch3"
$ws.Range("D10").Value = "ECs"
$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 29.31506333333334
$ws.Range("H10").Value = 87.94519000000001
$ws.Range("I10").Value = 0.2873975381543141
$ws.Range("J10").Value = 0.2873975381543141
$ws.Range("K10").Value = 2
$ws.Range("L10").Value = 0.6666666666666666
$ws.Range("M10").Value = 3.684362666666666
$ws.Range("N10").Value = 11.053088
$ws.Range("O10").Value = 0.0397011572965827
$ws.Range("P10").Value = 0.03970115729658269
$ws.Range("Q10").Value = 108.0073249163022
$ws.Range("R10").Value = 972.06592424672
$ws.Range("S10").Value = 0.01141001486891505
$ws.Range("T10").Value = 0.01141001486891505

# Row 11
$ws.Range("A11").Value = "M2"
$ws.Range("B11").Value = "Psen1"
$ws.Range("C11").Value = "Notch3"
$ws.Range("D11").Value = "FAPs"
$ws.Range("E11").Value = 3
$ws.Range("F11").Value = 1
$ws.Range("G11").Value = 29.31506333333334
$ws.Range("H11").Value = 87.94519000000001
$ws.Range("I11").Value = 0.2873975381543141
$ws.Range("J11").Value = 0.2873975381543141
$ws.Range("K11").Value = 3
$ws.Range("L11").Value = 1
$ws.Range("M11").Value = 5.038243666666667
$ws.Range("N11").Value = 15.114731
$ws.Range("O11").Value = 0.05429001496473517
$ws.Range("P11").Value = 0.05429001496473517
$ws.Range("Q11").Value = 147.6964321770989
$ws.Range("R11").Value = 1329.26788959389
$ws.Range("S11").Value = 0.01560281664722576
$ws.Range("T11").Value = 0.01560281664722576

# Row 12
$ws.Range("A12").Value = "M2"
$ws.Range("B12").Value = "Psen1"
$ws.Range("C12").Value = "Notch3"
$ws.Range("D12").Value = "M2"
$ws.Range("E12").Value = 3
$ws.Range("F12").Value = 1
$ws.Range("G12").Value = 29.31506333333334
$ws.Range("H12").Value = 87.94519000000001
$ws.Range("I12").Value = 0.2873975381543141
$ws.Range("J12").Value = 0.2873975381543141
$ws.Range("K12").Value = 2
$ws.Range("L12").Value = 0.6666666666666666
$ws.Range("M12").Value = 0.09716666666666667
$ws.Range("N12").Value = 0.2915
$ws.Range("O12").Value = 0.001047027523164011
$ws.Range("P12").Value = 0.001047027523164011
$ws.Range("Q12").Value = 2.848446987222223
$ws.Range("R12").Value = 25.636022885
$ws.Range("S12").Value = 0.000300913132537146
$ws.Range("T12").Value = 0.000300913132537146

# Row 13
$ws.Range("A13").Value = "M2"
$ws.Range("B13").Value = "Psen1"
$ws.Range("C13").Value = "Notch3"
$ws.Range("D13").Value = "sCs"
$ws.Range("E13").Value = 3
$ws.Range("F13").Value = 1
$ws.Range("G13").Value = 29.31506333333334
$ws.Range("H13").Value = 87.94519000000001
$ws.Range("I13").Value = 0.2873975381543141
$ws.Range("J13").Value = 0.2873975381543141
$ws.Range("K13").Value = 3
$ws.Range("L13").Value = 1
$ws.Range("M13").Value = 83.98262666666666
$ws.Range("N13").Value = 251.94788
$ws.Range("O13").Value = 0.9049618002155182
$ws.Range("P13").Value = 0.9049618002155182
$ws.Range("Q13").Value = 2461.956019633023
$ws.Range("R13").Value = 22157.6041766972
$ws.Range("S13").Value = 0.2600837935056362
$ws.Range("T13").Value = 0.2600837935056362

# Row 14
$ws.Range("A14").Value = "sCs"
$ws.Range("B14").Value = "Psen1"
$ws.Range("C14").Value = "Notch3"
$ws.Range("D14").Value = "ECs"
$ws.Range("E14").Value = 3
$ws.Range("F14").Value = 1
$ws.Range("G14").Value = 18.86569966666666
$ws.Range("H14").Value = 56.59709899999999
$ws.Range("I14").Value = 0.1849545940974826
$ws.Range("J14").Value = 0.1849545940974826
$ws.Range("K14").Value = 2
$ws.Range("L14").Value = 0.6666666666666666
$ws.Range("M14").Value = 3.684362666666666
$ws.Range("N14").Value = 11.053088
$ws.Range("O14").Value = 0.0397011572965827
$ws.Range("P14").Value = 0.03970115729658269
$ws.Range("Q14").Value = 69.50807953241244
$ws.Range("R14").Value = 625.5727157917119
$ws.Range("S14").Value = 0.007342911432989764
$ws.Range("T14").Value = 0.007342911432989763

# Row 15
$ws.Range("A15").Value = "sCs"
$ws.Range("B15").Value = "Psen1"
$ws.Range("C15").Value = "Notch3"
$ws.Range("D15").Value = "FAPs"
$ws.Range("E15").Value = 3
$ws.Range("F15").Value = 1
$ws.Range("G15").Value = 18.86569966666666
$ws.Range("H15").Value = 56.59709899999999
$ws.Range("I15").Value = 0.1849545940974826
$ws.Range("J15").Value = 0.1849545940974826
$ws.Range("K15").Value = 3
$ws.Range("L15").Value = 1
$ws.Range("M15").Value = 5.038243666666667
$ws.Range("N15").Value = 15.114731
$ws.Range("O15").Value = 0.05429001496473517
$ws.Range("P15").Value = 0.05429001496473517
$ws.Range("Q15").Value = 95.04999186281877
$ws.Range("R15").Value = 855.449926765369
$ws.Range("S15").Value = 0.01004118768134885
$ws.Range("T15").Value = 0.01004118768134885

# Row 16
$ws.Range("A16").Value = "sCs"
$ws.Range("B16").Value = "Psen1"
$ws.Range("C16").Value = "Notch3"
$ws.Range("D16").Value = "M2"
$ws.Range("E16").Value = 3
$ws.Range("F16").Value = 1
$ws.Range("G16").Value = 18.86569966666666
$ws.Range("H16").Value = 56.59709899999999
$ws.Range("I16").Value = 0.1849545940974826
$ws.Range("J16").Value = 0.1849545940974826
$ws.Range("K16").Value = 2
$ws.Range("L16").Value = 0.6666666666666666
$ws.Range("M16").Value = 0.09716666666666667
$ws.Range("N16").Value = 0.2915
$ws.Range("O16").Value = 0.001047027523164011
$ws.Range("P16").Value = 0.001047027523164011
$ws.Range("Q16").Value = 1.833117150944444
$ws.Range("R16").Value = 16.4980543585
$ws.Range("S16").Value = 0.0001936525505556923
$ws.Range("T16").Value = 0.0001936525505556923

# Row 17
$ws.Range("A17").Value = "sCs"
$ws.Range("B17").Value = "Psen1"
$ws.Range("C17").Value = "Notch3"
$ws.Range("D17").Value = "sCs"
$ws.Range("E17").Value = 3
$ws.Range("F17").Value = 1
$ws.Range("G17").Value = 18.86569966666666
$ws.Range("H17").Value = 56.59709899999999
$ws.Range("I17").Value = 0.1849545940974826
$ws.Range("J17").Value = 0.1849545940974826
$ws.Range("K17").Value = 3
$ws.Range("L17").Value = 1
$ws.Range("M17").Value = 83.98262666666666
$ws.Range("N17").Value = 251.94788
$ws.Range("O17").Value = 0.9049618002155182
$ws.Range("P17").Value = 0.9049618002155182
$ws.Range("Q17").Value = 1584.391011911124
$ws.Range("R17").Value = 14259.51910720012
$ws.Range("S17").Value = 0.1673768424325884
$ws.Range("T17").Value = 0.1673768424325884
